$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @('000072', 'Details: IMU CJMCU-20948 Raw Data Reading - Raw Plots of 100 points of Accel Gyro and Magnometer data separately. Script used: Read_IMU.  Dataset used: Arduino Serial Output of IMU CJMCU-20948. File Location: Visualisations/IMU_RealRawData. Date Generated: 03-Mar-2023 10:17:11'),
    @('000073', 'Details: IMU CJMCU-20948 Raw Data Reading - Raw Plots of 100 points of Accel Gyro and Magnometer data separately. Script used: Read_IMU.  Dataset used: Arduino Serial Output of IMU CJMCU-20948. File Location: Visualisations/IMU_RealRawData. Date Generated: 03-Mar-2023 10:17:30'),
    @('000074', 'Details: IMU CJMCU-20948 Raw Data Reading - Raw Plots of 100 points of Accel Gyro and Magnometer data separately. Script used: Read_IMU.  Dataset used: Arduino Serial Output of IMU CJMCU-20948. File Location: Visualisations/IMU_RealRawData. Date Generated: 03-Mar-2023 10:18:34'),
    @('000075', 'Details: IMU CJMCU-20948 Raw Data Reading - Raw Plots of 100 points of Accel Gyro and Magnometer data separately. Script used: Read_IMU.  Dataset used: Arduino Serial Output of IMU CJMCU-20948. File Location: Visualisations/IMU_RealRawData. Date Generated: 03-Mar-2023 10:18:48'),
    @('000076', 'Details: IMU CJMCU-20948 Raw Data Reading - Raw Plots of 100 points of Accel Gyro and Magnometer data separately. Script used: Read_IMU.  Dataset used: Arduino Serial Output of IMU CJMCU-20948. File Location: Visualisations/IMU_RealRawData. Date Generated: 03-Mar-2023 10:19:18'),
    @('000077', 'Details: IMU CJMCU-20948 Raw Data Reading - Raw Plots of 100 points of Accel Gyro and Magnometer data separately. Script used: Read_IMU.  Dataset used: Arduino Serial Output of IMU CJMCU-20948. File Location: Visualisations/IMU_RealRawData. Date Generated: 03-Mar-2023 10:19:45'),
    @('000078', 'Details: IMU CJMCU-20948 Raw Data Reading - Raw Plots of 100 points of Accel Gyro and Magnometer data separately. Script used: Read_IMU.  Dataset used: Arduino Serial Output of IMU CJMCU-20948. File Location: Visualisations/IMU_RealRawData. Date Generated: 03-Mar-2023 10:22:40'),
    @('000079', 'Details: IMU CJMCU-20948 Raw Data Reading - Raw Plots of 100 points of Accel Gyro and Magnometer data separately. Script used: Read_IMU.  Dataset used: Arduino Serial Output of IMU CJMCU-20948. File Location: Visualisations/IMU_RealRawData. Date Generated: 03-Mar-2023 10:22:58'),
    @('000080', 'Details: IMU CJMCU-20948 Raw Data Reading - Raw Plots of 100 points of Accel Gyro and Magnometer data separately. Script used: Read_IMU.  Dataset used: Arduino Serial Output of IMU CJMCU-20948. File Location: Visualisations/IMU_RealRawData. Date Generated: 03-Mar-2023 10:23:32'),
    @('000081', 'Details: IMU CJMCU-20948 Raw Data Reading - Raw Plots of 100 points of Accel Gyro and Magnometer data separately. Script used: Read_IMU.  Dataset used: Arduino Serial Output of IMU CJMCU-20948. File Location: Visualisations/IMU_RealRawData. Date Generated: 03-Mar-2023 10:23:57'),
    @('000082', 'Details: IMU CJMCU-20948 Raw Data Reading - Raw Plots of 100 points of Accel Gyro and Magnometer data separately. Script used: Read_IMU.  Dataset used: Arduino Serial Output of IMU CJMCU-20948. File Location: Visualisations/IMU_RealRawData. Date Generated: 03-Mar-2023 10:26:40')
)

$startRow = 73
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $data[$i][0]
    $cellA.Style = "Normal"
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
